# Auto-generated edit script applying the Pandaemonium_Profits diff
# across all 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 681
$ws.Range("I28").Value = 442.6111
$ws.Range("K28").Value = 442.6111
$ws.Range("M28").Value = 42.38889999999998
$ws.Range("H32").Value = 847.5
$ws.Range("I32").Value = 833.3333
$ws.Range("J32").Value = 856
$ws.Range("K32").Value = 833.3333
$ws.Range("L32").Value = 856
$ws.Range("M32").Value = -507.3333
$ws.Range("N32").Value = -1508
$ws.Range("H33").Value = 370.65216
$ws.Range("I33").Value = 325.66666
$ws.Range("K33").Value = 325.66666
$ws.Range("M33").Value = -96.66665999999998
$ws.Range("H64").Value = 4091.3914
$ws.Range("I64").Value = 3585.8572
$ws.Range("J64").Value = 4877.778
$ws.Range("K64").Value = 3585.8572
$ws.Range("L64").Value = 4877.778
$ws.Range("M64").Value = -3337.8572
$ws.Range("N64").Value = -5373.778
$ws.Range("H67").Value = 4091.3914
$ws.Range("I67").Value = 3585.8572
$ws.Range("J67").Value = 4877.778
$ws.Range("K67").Value = 3585.8572
$ws.Range("L67").Value = 4877.778
$ws.Range("M67").Value = -2727.8572
$ws.Range("N67").Value = -6593.778
$ws.Range("H98").Value = 2358.4
$ws.Range("I98").Value = 1842.6666
$ws.Range("J98").Value = 7000
$ws.Range("K98").Value = 1842.6666
$ws.Range("L98").Value = 7000
$ws.Range("M98").Value = -344.6666
$ws.Range("N98").Value = -9996
$ws.Range("H113").Value = 2788.7144
$ws.Range("I113").Value = 2836.6667
$ws.Range("K113").Value = 2836.6667
$ws.Range("M113").Value = 417.3332999999998
$ws.Range("H122").Value = 2358.4
$ws.Range("I122").Value = 1842.6666
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 5527.9998
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -3077.9998
$ws.Range("N122").Value = -25900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 10000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 10000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 10000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -10288
$ws.Range("H32").Value = 7238.9507
$ws.Range("I32").Value = 5618.067
$ws.Range("J32").Value = 27500
$ws.Range("K32").Value = 5618.067
$ws.Range("L32").Value = 27500
$ws.Range("M32").Value = -5331.067
$ws.Range("N32").Value = -28074
$ws.Range("H45").Value = 1201.8
$ws.Range("I45").Value = 1103.1666
$ws.Range("J45").Value = 1596.3334
$ws.Range("K45").Value = 1103.1666
$ws.Range("L45").Value = 1596.3334
$ws.Range("M45").Value = -726.1666
$ws.Range("N45").Value = -2350.3334
$ws.Range("H63").Value = 3333.8333
$ws.Range("J63").Value = 4000.75
$ws.Range("L63").Value = 4000.75
$ws.Range("N63").Value = -5372.75
$ws.Range("H66").Value = 3333.8333
$ws.Range("J66").Value = 4000.75
$ws.Range("L66").Value = 20003.75
$ws.Range("N66").Value = -26867.75
$ws.Range("H80").Value = 142882560
$ws.Range("I80").Value = 18750
$ws.Range("K80").Value = 18750
$ws.Range("M80").Value = -17752
$ws.Range("H83").Value = 142882560
$ws.Range("I83").Value = 18750
$ws.Range("K83").Value = 56250
$ws.Range("M83").Value = -51258
$ws.Range("H97").Value = 890.86365
$ws.Range("I97").Value = 715.7368
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 715.7368
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -219.7368
$ws.Range("N97").Value = -2992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H120").Value = 50761
$ws.Range("J120").Value = 50761
$ws.Range("L120").Value = 50761
$ws.Range("N120").Value = -60437

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3249.25
$ws.Range("I2").Value = 999
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 999
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = -886
$ws.Range("N2").Value = -10226
$ws.Range("H122").Value = 12431.158
$ws.Range("I122").Value = 7048.625
$ws.Range("J122").Value = 16345.728
$ws.Range("K122").Value = 21145.875
$ws.Range("L122").Value = 49037.18399999999
$ws.Range("M122").Value = -18695.875
$ws.Range("N122").Value = -53937.18399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 70000
$ws.Range("J37").Value = 70000
$ws.Range("L37").Value = 210000
$ws.Range("N37").Value = -210224
$ws.Range("H38").Value = 75.888885
$ws.Range("I38").Value = 29.375
$ws.Range("J38").Value = 113.1
$ws.Range("K38").Value = 88.125
$ws.Range("L38").Value = 339.3
$ws.Range("M38").Value = 258.875
$ws.Range("N38").Value = -1033.3
$ws.Range("H121").Value = 1933.3334
$ws.Range("I121").Value = 1400
$ws.Range("J121").Value = 2000
$ws.Range("K121").Value = 4200
$ws.Range("L121").Value = 6000
$ws.Range("M121").Value = -2890
$ws.Range("N121").Value = -8620
$ws.Range("H131").Value = 18159.629
$ws.Range("I131").Value = 396.27658
$ws.Range("J131").Value = 137427.86
$ws.Range("K131").Value = 1188.82974
$ws.Range("L131").Value = 412283.58
$ws.Range("M131").Value = 3851.17026
$ws.Range("N131").Value = -422363.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11715
$ws.Range("I80").Value = 18001.666
$ws.Range("J80").Value = 7000
$ws.Range("K80").Value = 18001.666
$ws.Range("L80").Value = 7000
$ws.Range("M80").Value = -17003.666
$ws.Range("N80").Value = -8996
$ws.Range("H83").Value = 11715
$ws.Range("I83").Value = 18001.666
$ws.Range("J83").Value = 7000
$ws.Range("K83").Value = 90008.33
$ws.Range("L83").Value = 35000
$ws.Range("M83").Value = -85016.33
$ws.Range("N83").Value = -44984
$ws.Range("H122").Value = 8357.700000000001
$ws.Range("I122").Value = 9653.857
$ws.Range("J122").Value = 5333.3335
$ws.Range("K122").Value = 28961.571
$ws.Range("L122").Value = 16000.0005
$ws.Range("M122").Value = -26511.571
$ws.Range("N122").Value = -20900.0005
$ws.Range("H132").Value = 5974.76
$ws.Range("I132").Value = 15967.857
$ws.Range("J132").Value = 2088.5557
$ws.Range("K132").Value = 47903.571
$ws.Range("L132").Value = 6265.6671
$ws.Range("M132").Value = -45373.571
$ws.Range("N132").Value = -11325.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3877.5
$ws.Range("I7").Value = 3500
$ws.Range("J7").Value = 4129.1665
$ws.Range("K7").Value = 3500
$ws.Range("L7").Value = 4129.1665
$ws.Range("M7").Value = -3388
$ws.Range("N7").Value = -4353.1665
$ws.Range("H19").Value = 10000
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H55").Value = 191158.14
$ws.Range("I55").Value = 334036.75
$ws.Range("J55").Value = 653.3333
$ws.Range("K55").Value = 334036.75
$ws.Range("L55").Value = 653.3333
$ws.Range("M55").Value = -333863.75
$ws.Range("N55").Value = -999.3333
$ws.Range("H68").Value = 2500
$ws.Range("J68").Value = 3000
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 2500
$ws.Range("J71").Value = 3000
$ws.Range("L71").Value = 15000
$ws.Range("N71").Value = -22488
$ws.Range("H93").Value = 2999.0908
$ws.Range("I93").Value = 3071.4285
$ws.Range("J93").Value = 2872.5
$ws.Range("K93").Value = 3071.4285
$ws.Range("L93").Value = 2872.5
$ws.Range("M93").Value = -1823.4285
$ws.Range("N93").Value = -5368.5
$ws.Range("H100").Value = 4382.278
$ws.Range("I100").Value = 2656.75
$ws.Range("K100").Value = 2656.75
$ws.Range("M100").Value = -2115.75
$ws.Range("H108").Value = 31000
$ws.Range("J108").Value = 31000
$ws.Range("L108").Value = 31000
$ws.Range("N108").Value = -38680
$ws.Range("H126").Value = 3877.5
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 4129.1665
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 12387.4995
$ws.Range("M126").Value = -8030
$ws.Range("N126").Value = -17327.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3500
$ws.Range("I81").Value = 3000
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 6000
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = -4939
$ws.Range("N81").Value = -10122
$ws.Range("H84").Value = 3500
$ws.Range("I84").Value = 3000
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 30000
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = -24696
$ws.Range("N84").Value = -50608
$ws.Range("H96").Value = 3434.3333
$ws.Range("I96").Value = 2651.5
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 2651.5
$ws.Range("L96").Value = 5000
$ws.Range("M96").Value = -1278.5
$ws.Range("N96").Value = -7746
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H126").Value = 1386.2858
$ws.Range("I126").Value = 1140.8
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3422.4
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -952.3999999999996
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 2724.075
$ws.Range("I132").Value = 2404.318
$ws.Range("J132").Value = 3114.889
$ws.Range("K132").Value = 7212.954000000001
$ws.Range("L132").Value = 9344.667000000001
$ws.Range("M132").Value = -4682.954000000001
$ws.Range("N132").Value = -14404.667
